$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manchester tribunal venue details
$ws.Range("B3").Value = "Manchester Employment Tribunal, Alexandra House, 14-22 The Parsonage, Manchester, M3 2JA"
$ws.Range("B4").Value = "0161 833 6100"
$ws.Range("B5").Value = "0870 739 4433"
$ws.Range("B6").Value = "DX 743570"
$ws.Range("B7").Value = "Manchesteret@justice.gov.uk"

# Glasgow tribunal venue details
$ws.Range("B8").Value = "Eagle Building, 215 Bothwell Street, Glasgow, G2 7TS"
$ws.Range("B9").Value = "0141 204 0730"
$ws.Range("B10").Value = "01264 785 177"
$ws.Range("B11").Value = "DX 7435701"
$ws.Range("B12").Value = "glasgowet@justice.gov.uk"

# Recreate the hyperlinks so the "display" text shown for the email
# addresses reflects the new values, while the underlying mailto:
# targets stay the same as before.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:manchester@gmail.com", "", "", "Manchesteret@justice.gov.uk")
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:glasgow@gmail.com", "", "", "glasgowet@justice.gov.uk")

# Adding a hyperlink forces Excel's built-in "Hyperlink" cell style onto
# the cell; restore the original (default) look of these two cells so
# formatting matches the rest of the sheet.
foreach ($addr in @("B7", "B12")) {
    $f = $ws.Range($addr).Font
    $f.Name = "Calibri"
    $f.Size = 11
    $f.Color = 0
    $f.Bold = $false
    $f.Italic = $false
    $f.Underline = -4142
}
